$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Sheet1")

# Updated Participants (CDS) query text - fixed primary diagnosis filter
$newQuery = @"
MATCH (p:participant)-->(s:study)
OPTIONAL MATCH (samp:sample)-->(p)
OPTIONAL MATCH (p)<--(diag:diagnosis)
OPTIONAL MATCH (samp)<--(f:file)
OPTIONAL MATCH (f)<--(g:genomic_info)
WITH s, p, samp, f, g, diag
WHERE diag.primary_diagnosis in ['Adrenal Cortical Carcinoma']
WITH p
OPTIONAL MATCH (p)-->(s:study)
OPTIONAL MATCH (samp:sample)-->(p)
WITH s, p, apoc.coll.sort(collect(distinct samp.sample_id)) as samp
RETURN 
coalesce(p.participant_id,'') as ``Participant ID``,
coalesce(s.study_name, '') as ``Study Name``,
coalesce(s.phs_accession,'') as ``Accession``,
coalesce(p.gender,'') as ``Gender``,
coalesce(apoc.text.join(samp, ','), '') as ``Samples``
ORDER BY p.participant_id
LIMIT 100
"@

$ws.Range("B2").Value = $newQuery

# Keep row 2 height as authored (unaffected by the longer query text)
$ws.Rows.Item(2).RowHeight = 157.5

# Active cell moved to E4 in the saved view
$ws.Range("E4").Select()
